# Update "想去人数" (want-to-go count) figures in column F across sheets
# "展览" (sheet 1), "演出" (sheet 2), and "全部类型" (sheet 4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1323
$ws1.Range("F4").Value  = 77
$ws1.Range("F8").Value  = 123
$ws1.Range("F9").Value  = 1002
$ws1.Range("F10").Value = 325
$ws1.Range("F17").Value = 134
$ws1.Range("F21").Value = 971
$ws1.Range("F22").Value = 431
$ws1.Range("F27").Value = 32
$ws1.Range("F28").Value = 453

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 33
$ws2.Range("F9").Value = 7

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1323
$ws4.Range("F6").Value  = 78
$ws4.Range("F10").Value = 123
$ws4.Range("F11").Value = 1002
$ws4.Range("F12").Value = 325
$ws4.Range("F18").Value = 33
$ws4.Range("F24").Value = 134
$ws4.Range("F28").Value = 971
$ws4.Range("F29").Value = 431
$ws4.Range("F31").Value = 7
$ws4.Range("F38").Value = 32
$ws4.Range("F40").Value = 453
